$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency market data to rows 2-51 (columns B-E).
# Columns D (Price) and E (Volume) are formatted as Text to preserve
# the exact original string formatting (e.g. "57.916.93", "  +2.36%  ").

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.916.93"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.36%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.066.65"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.93%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "526.69"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +6.12%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.54"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +6.14%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.00%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +5.68%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.66"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +6.22%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +7.41%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.371"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.96%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.09%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.03%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.39"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +8.71%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000171"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +16.64%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "57.927.73"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.08%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.22"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +7.16%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.069.70"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.07%  "

# Row 19
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +7.34%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.21"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.20%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "341.64"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.71%  "

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.13%  "

# Row 23
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.504"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +7.29%  "

# Row 24
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.95"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +5.44%  "

# Row 25
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0₃0980"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +9.36%  "

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +6.06%  "

# Row 27
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.25%  "

# Row 28
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.99"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +7.41%  "

# Row 29
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.44"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +9.92%  "

# Row 30
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.86"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +6.83%  "

# Row 31
$ws.Range("B31").Value = "Fetch.AI"
$ws.Range("C31").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.24"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +6.76%  "

# Row 32
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.11"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.02%  "

# Row 33
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.80"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +7.68%  "

# Row 34
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "156.93"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.96%  "

# Row 35
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.01"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +7.28%  "

# Row 36
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.33"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.06%  "

# Row 37
$ws.Range("B37").Value = "EnergySwap"
$ws.Range("C37").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "26.40"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +13.52%  "

# Row 38
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0707"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +5.31%  "

# Row 39
$ws.Range("B39").Value = "RenzoRestakedETH"
$ws.Range("C39").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.103.25"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.13%  "

# Row 40
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.88"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.25%  "

# Row 41
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.93"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +10.31%  "

# Row 42
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.668"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.28%  "

# Row 43
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.09%  "

# Row 44
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.48"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +5.58%  "

# Row 45
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.340.39"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.32%  "

# Row 46
$ws.Range("B46").Value = "ONDO"
$ws.Range("C46").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.03"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.26%  "

# Row 47
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.01"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.09%  "

# Row 48
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.05"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.39%  "

# Row 49
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0245"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.91%  "

# Row 50
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.21"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +6.37%  "

# Row 51
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0901"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +6.14%  "

